# Natmi following Dr Hou advice
#
# The original sheet only modelled two "Sending cluster" groups (ECs, sCs)
# against three "Target cluster" groups (ECs, FAPs, sCs) for the Cd274->Cd80
# ligand-receptor pair (rows 2-7). Dr Hou's advice was to also treat "FAPs"
# as a Sending cluster, so the table becomes a full 3x3 cross of
# Sending cluster x Target cluster (ECs, FAPs, sCs) - i.e. rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ligand = "Cd274"
$receptor = "Cd80"

# Each entry: Sending cluster, Target cluster, then the 16 numeric metric
# columns (E..T) for that Sending/Target pair.
$rows = @(
    @("ECs", "ECs", 3, 1, 32.09065066666667, 96.271952, 0.946743898370603, 0.9467438983706029, 3, 1, 26.499453, 79.49835900000001, 0.7877954840311897, 0.7877954840311898, 850.3846890807521, 7653.462201726768, 0.7458405676704447, 0.7458405676704447),
    @("ECs", "FAPs", 3, 1, 32.09065066666667, 96.271952, 0.946743898370603, 0.9467438983706029, 3, 1, 5.993807666666666, 17.981423, 0.1781883804149286, 0.1781883804149287, 192.3451879941884, 1731.106691947696, 0.1686987619183735, 0.1686987619183735),
    @("ECs", "sCs", 3, 1, 32.09065066666667, 96.271952, 0.946743898370603, 0.9467438983706029, 3, 1, 1.144217, 3.432650999999999, 0.03401613555388164, 0.03401613555388164, 36.71866803386133, 330.4680123047519, 0.03220456878178478, 0.03220456878178478),
    @("FAPs", "ECs", 2, 0.6666666666666666, 1.154570666666667, 3.463712, 0.03406234249527876, 0.03406234249527877, 3, 1, 26.499453, 79.49835900000001, 0.7877954840311897, 0.7877954840311898, 30.595491116512, 275.3594200486081, 0.0268341595933043, 0.0268341595933043),
    @("FAPs", "FAPs", 2, 0.6666666666666666, 1.154570666666667, 3.463712, 0.03406234249527876, 0.03406234249527877, 3, 1, 5.993807666666666, 17.981423, 0.1781883804149286, 0.1781883804149287, 6.92027451357511, 62.282470622176, 0.006069513642372321, 0.006069513642372324),
    @("FAPs", "sCs", 2, 0.6666666666666666, 1.154570666666667, 3.463712, 0.03406234249527876, 0.03406234249527877, 3, 1, 1.144217, 3.432650999999999, 0.03401613555388164, 0.03401613555388164, 1.321079384501333, 11.889714460512, 0.001158669259602146, 0.001158669259602146),
    @("sCs", "ECs", 3, 1, 0.6505880000000001, 1.951764, 0.01919375913411833, 0.01919375913411833, 3, 1, 26.499453, 79.49835900000001, 0.7877954840311897, 0.7877954840311898, 17.240226128364, 155.162035155276, 0.01512075676744082, 0.01512075676744082),
    @("sCs", "FAPs", 3, 1, 0.6505880000000001, 1.951764, 0.01919375913411833, 0.01919375913411833, 3, 1, 5.993807666666666, 17.981423, 0.1781883804149286, 0.1781883804149287, 3.899499342241334, 35.095494080172, 0.003420104854182789, 0.003420104854182789),
    @("sCs", "sCs", 3, 1, 0.6505880000000001, 1.951764, 0.01919375913411833, 0.01919375913411833, 3, 1, 1.144217, 3.432650999999999, 0.03401613555388164, 0.03401613555388164, 0.744413849596, 6.699724646363999, 0.0006528975124947231, 0.0006528975124947231)
)

$excelRow = 2
foreach ($entry in $rows) {
    $sendingCluster = $entry[0]
    $targetCluster = $entry[1]

    $ws.Cells.Item($excelRow, 1).Value = $sendingCluster
    $ws.Cells.Item($excelRow, 2).Value = $ligand
    $ws.Cells.Item($excelRow, 3).Value = $receptor
    $ws.Cells.Item($excelRow, 4).Value = $targetCluster

    # Columns E..T (5..20) hold the 16 numeric metrics that follow the
    # cluster/gene columns in each data row.
    for ($i = 0; $i -lt 16; $i++) {
        $col = 5 + $i
        $ws.Cells.Item($excelRow, $col).Value = $entry[2 + $i]
    }

    $excelRow++
}
